$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for LIDC-IDRI ...189/055.xml (row 26), shifting rows below it up.
$ws.Rows(26).Delete()

# Replace comma-space separators with " | " in column B (Rterms), rows 2..188.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace ', ', ' | '
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
